$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "'293.86"
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,5).Value = "'1.07%"
$ws.Cells.Item(2,5).Style = "Normal"
$ws.Cells.Item(3,4).Value = "'40.16"
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).Value = "'1.55%"
$ws.Cells.Item(3,5).Style = "Normal"
$ws.Cells.Item(4,4).Value = "'5.012"
$ws.Cells.Item(4,4).Style = "Normal"
$ws.Cells.Item(4,5).Value = "'-0.13%"
$ws.Cells.Item(4,5).Style = "Normal"
$ws.Cells.Item(5,4).Value = "'0.07312"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value = "'-0.49%"
$ws.Cells.Item(5,5).Style = "Normal"
$ws.Cells.Item(6,4).Value = "'4.319"
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Value = "'0.76%"
$ws.Cells.Item(6,5).Style = "Normal"
$ws.Cells.Item(7,4).Value = "'1.541"
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).Value = "'-1.03%"
$ws.Cells.Item(7,5).Style = "Normal"
$ws.Cells.Item(8,4).Value = "'0.9217"
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(10,4).Value = "'0.1169"
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).Value = "'-1.67%"
$ws.Cells.Item(10,5).Style = "Normal"
$ws.Cells.Item(11,4).Value = "'0.1789"
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Value = "'3.37%"
$ws.Cells.Item(11,5).Style = "Normal"
$ws.Cells.Item(12,4).Value = "'0.08667"
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).Value = "'-0.03%"
$ws.Cells.Item(12,5).Style = "Normal"
$ws.Cells.Item(13,4).Value = "'0.04306"
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,5).Value = "'3.14%"
$ws.Cells.Item(13,5).Style = "Normal"
$ws.Cells.Item(14,4).Value = "'0.1052"
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).Value = "'0.09%"
$ws.Cells.Item(14,5).Style = "Normal"
$ws.Cells.Item(15,4).Value = "'0.001284"
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).Value = "'0.57%"
$ws.Cells.Item(15,5).Style = "Normal"
$ws.Cells.Item(16,4).Value = "'0.005932"
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).Value = "'1.61%"
$ws.Cells.Item(16,5).Style = "Normal"
$ws.Cells.Item(17,2).Value = "HotbitToken"
$ws.Cells.Item(17,3).Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Cells.Item(17,4).Value = "'0.003799"
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,5).Value = "'3.08%"
$ws.Cells.Item(17,5).Style = "Normal"
$ws.Cells.Item(18,2).Value = "LEO"
$ws.Cells.Item(18,3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(18,4).Value = "'3.340"
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,5).Value = "'-1.28%"
$ws.Cells.Item(18,5).Style = "Normal"
$ws.Cells.Item(19,2).Value = "BitpandaEcosystemToken"
$ws.Cells.Item(19,3).Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Cells.Item(19,4).Value = "'0.3292"
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).Value = "'-0.16%"
$ws.Cells.Item(19,5).Style = "Normal"
$ws.Cells.Item(20,2).Value = "MCDex"
$ws.Cells.Item(20,3).Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Cells.Item(20,4).Value = "'7.926"
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).Value = "'5.19%"
$ws.Cells.Item(20,5).Style = "Normal"
$ws.Cells.Item(21,2).Value = "ProBitToken"
$ws.Cells.Item(21,3).Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Cells.Item(21,4).Value = "'0.1385"
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).Value = "'2.29%"
$ws.Cells.Item(21,5).Style = "Normal"
$ws.Cells.Item(22,2).Value = "ZBToken"
$ws.Cells.Item(22,3).Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Cells.Item(22,4).Value = "'0.2809"
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).Value = "'2.67%"
$ws.Cells.Item(22,5).Style = "Normal"
$ws.Cells.Item(23,2).Value = "CoinExToken"
$ws.Cells.Item(23,3).Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Cells.Item(23,4).Value = "'0.03961"
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).Value = "'3.06%"
$ws.Cells.Item(23,5).Style = "Normal"
$ws.Cells.Item(24,2).Value = "BitKan"
$ws.Cells.Item(24,3).Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Cells.Item(24,4).Value = "'0.001271"
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).Value = "'-0.21%"
$ws.Cells.Item(24,5).Style = "Normal"
$ws.Cells.Item(25,4).Value = "'0.0001232"
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).Value = "'-3.90%"
$ws.Cells.Item(25,5).Style = "Normal"
$ws.Cells.Item(26,4).Value = "'0.0003737"
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).Value = "'0.12%"
$ws.Cells.Item(26,5).Style = "Normal"
$ws.Cells.Item(38,4).Value = "'0.02325"
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Cells.Item(38,5).Value = "'0.63%"
$ws.Cells.Item(38,5).Style = "Normal"
$ws.Cells.Item(39,4).Value = "'0.05063"
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).Value = "'1.02%"
$ws.Cells.Item(39,5).Style = "Normal"
$ws.Cells.Item(40,4).Value = "'0.005987"
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).Value = "'17.17%"
$ws.Cells.Item(40,5).Style = "Normal"
$ws.Cells.Item(41,4).Value = "'0.007759"
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).Value = "'0.71%"
$ws.Cells.Item(41,5).Style = "Normal"
$ws.Cells.Item(42,4).Value = "'0.1290"
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,5).Value = "'1.71%"
$ws.Cells.Item(42,5).Style = "Normal"
$ws.Cells.Item(43,4).Value = "'0.007409"
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(44,4).Value = "'0.006958"
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).Value = "'-9.31%"
$ws.Cells.Item(44,5).Style = "Normal"
$ws.Cells.Item(45,4).Value = "'0.2920"
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).Value = "'-7.65%"
$ws.Cells.Item(45,5).Style = "Normal"
$ws.Cells.Item(46,4).Value = "'0.00006216"
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).Value = "'-5.03%"
$ws.Cells.Item(46,5).Style = "Normal"
$ws.Cells.Item(47,4).Value = "'0.00000000753"
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,5).Value = "'0.28%"
$ws.Cells.Item(47,5).Style = "Normal"
$ws.Cells.Item(48,4).Value = "'0.04639"
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(48,5).Value = "'-81.58%"
$ws.Cells.Item(48,5).Style = "Normal"
$ws.Cells.Item(49,4).Value = "'0.004216"
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).Value = "'0.13%"
$ws.Cells.Item(49,5).Style = "Normal"
$ws.Cells.Item(50,4).Value = "'0.00002109"
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).Value = "'0.28%"
$ws.Cells.Item(50,5).Style = "Normal"
$ws.Cells.Item(51,4).Value = "'0.0002008"
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,5).Value = "'0.28%"
$ws.Cells.Item(51,5).Style = "Normal"
